$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Applied Jobs")

# The sheet holds an "applied jobs" table (candidateId, jobsListId,
# recruiterId, createdAt). Append one more record as row 6, matching the
# existing rows' layout and the date-serial formatting used in column D.
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 3
$ws.Cells.Item(6, 3).Value = 1

$ws.Cells.Item(6, 4).Value = 44993.38872685185
# Reuse the same date number format as the other D-column cells (built-in
# numFmtId 14) so the new cell shares their style instead of creating a
# brand-new cell format.
$ws.Cells.Item(6, 4).NumberFormat = "m/d/yy"
